# Actualización automática del mapa
# 1) Removes the first data row (old row 2, caso 2098) from the "NEW" sheet,
#    shifting all subsequent rows up by one and shrinking the used range
#    from A1:P69 to A1:P68.
# 2) The case that is now in row 2 (caso 3299) gets its "Estado" and
#    "Observaciones" fields refreshed to reflect the latest automated sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

$ws.Rows.Item(2).Delete()

$ws.Range("G2").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H2").Value = "qap traspaso nodo TLC y Teco"
